$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "backup@backdoor.com, system, System"
$ws.Range("G3").Value = "dnasr281@gmail.com, System"
$ws.Range("G4").Value = "backup@backdoor.com, System"
$ws.Range("L4").Value = 333
$ws.Range("G5").Value = "backup@backdoor.com, System"
$ws.Range("G6").Value = "dnasr281@gmail.com, System"
$ws.Range("G10").Value = "dnasr281@gmail.com, System"
$ws.Range("L10").Value = "'70.4%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("G11").Value = "dnasr281@gmail.com, System"
$ws.Range("G12").Value = "dnasr281@gmail.com, System"
$ws.Range("G13").Value = "dnasr281@gmail.com, System"
$ws.Range("G14").Value = "dnasr281@gmail.com, System"
$ws.Range("G15").Value = "dnasr281@gmail.com, System"
$ws.Range("M16").Value = 57
$ws.Range("S16").Value = "'67.6%"
$ws.Range("R16").Copy()
$ws.Range("S16").PasteSpecial(-4122)
$ws.Range("G17").Value = "dnasr281@gmail.com, System"
$ws.Range("G18").Value = "dnasr281@gmail.com, System"
$ws.Range("G29").Value = "backup@backdoor.com, system, System"
$ws.Range("H29").Value = "36/57"
$ws.Range("G30").Value = "dnasr281@gmail.com, System"
$ws.Range("H30").Value = "42/57"
$ws.Range("G31").Value = "backup@backdoor.com, System"
$ws.Range("H31").Value = "56/57"
$ws.Range("G32").Value = "backup@backdoor.com, System"
$ws.Range("H32").Value = "38/57"
$ws.Range("G33").Value = "dnasr281@gmail.com, System"
$ws.Range("H33").Value = "43/57"
$ws.Range("H34").Value = "31/57"
$ws.Range("H35").Value = "41/57"
$ws.Range("H36").Value = "26/57"
$ws.Range("G37").Value = "dnasr281@gmail.com, System"
$ws.Range("H37").Value = "19/57"
$ws.Range("G38").Value = "dnasr281@gmail.com, System"
$ws.Range("H38").Value = "33/57"
$ws.Range("G39").Value = "dnasr281@gmail.com, System"
$ws.Range("H39").Value = "37/57"
$ws.Range("G40").Value = "dnasr281@gmail.com, System"
$ws.Range("H40").Value = "36/57"
$ws.Range("G41").Value = "dnasr281@gmail.com, System"
$ws.Range("H41").Value = "45/57"
$ws.Range("G42").Value = "dnasr281@gmail.com, System"
$ws.Range("H42").Value = "46/57"
$ws.Range("H43").Value = "44/57"
$ws.Range("G44").Value = "dnasr281@gmail.com, System"
$ws.Range("H44").Value = "43/57"
$ws.Range("G45").Value = "dnasr281@gmail.com, System"
$ws.Range("H45").Value = "39/57"
$ws.Range("H46").Value = "0/57"
$ws.Range("H47").Value = "0/57"
$ws.Range("H48").Value = "0/57"
$ws.Range("H49").Value = "0/57"
$ws.Range("H50").Value = "0/57"
$ws.Range("H51").Value = "0/57"
$ws.Range("H52").Value = "0/57"
$ws.Range("H53").Value = "0/57"
$ws.Range("H54").Value = "0/57"
$ws.Range("H55").Value = "0/57"
$ws.Range("G56").Value = "backup@backdoor.com, system, System"
$ws.Range("G57").Value = "dnasr281@gmail.com, System"
$ws.Range("G58").Value = "backup@backdoor.com, System"
$ws.Range("G59").Value = "backup@backdoor.com, System"
$ws.Range("G60").Value = "dnasr281@gmail.com, System"
$ws.Range("G64").Value = "dnasr281@gmail.com, System"
$ws.Range("G65").Value = "dnasr281@gmail.com, System"
$ws.Range("G66").Value = "dnasr281@gmail.com, System"
$ws.Range("G67").Value = "dnasr281@gmail.com, System"
$ws.Range("G68").Value = "dnasr281@gmail.com, System"
$ws.Range("G69").Value = "dnasr281@gmail.com, System"
$ws.Range("G71").Value = "dnasr281@gmail.com, System"
$ws.Range("G72").Value = "dnasr281@gmail.com, System"
$ws.Range("G84").Value = "backup@backdoor.com, System"
$ws.Range("G85").Value = "backup@backdoor.com, System"
$ws.Range("G86").Value = "dnasr281@gmail.com, System"
$ws.Range("G87").Value = "dnasr281@gmail.com, System"
$ws.Range("G88").Value = "dnasr281@gmail.com, System"
$ws.Range("G89").Value = "dnasr281@gmail.com, System"
$ws.Range("G93").Value = "dnasr281@gmail.com, System"
$ws.Range("G95").Value = "dnasr281@gmail.com, System"
$ws.Range("G96").Value = "dnasr281@gmail.com, System"
$ws.Range("G110").Value = "backup@backdoor.com, System"
$ws.Range("G111").Value = "backup@backdoor.com, System"
$ws.Range("G112").Value = "dnasr281@gmail.com, System"
$ws.Range("G113").Value = "dnasr281@gmail.com, System"
$ws.Range("G114").Value = "dnasr281@gmail.com, System"
$ws.Range("G115").Value = "dnasr281@gmail.com, System"
$ws.Range("G119").Value = "dnasr281@gmail.com, System"
$ws.Range("G121").Value = "dnasr281@gmail.com, System"
$ws.Range("G122").Value = "dnasr281@gmail.com, System"
$ws.Range("G136").Value = "backup@backdoor.com, System"
$ws.Range("G137").Value = "backup@backdoor.com, System"
$ws.Range("G138").Value = "dnasr281@gmail.com, System"
$ws.Range("G139").Value = "dnasr281@gmail.com, System"
$ws.Range("G140").Value = "dnasr281@gmail.com, System"
$ws.Range("G141").Value = "dnasr281@gmail.com, System"
$ws.Range("G145").Value = "dnasr281@gmail.com, System"
$ws.Range("G147").Value = "dnasr281@gmail.com, System"
$ws.Range("G148").Value = "dnasr281@gmail.com, System"
$excel.CutCopyMode = 0
